$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) — sheetId 1
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 11849
$ws1.Range("F5").Value = 291
$ws1.Range("F6").Value = 601
$ws1.Range("F8").Value = 299
$ws1.Range("F9").Value = 1067

# Sheet "全部类型" (All types) — sheetId 4
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 11849
$ws4.Range("F8").Value = 291
$ws4.Range("F9").Value = 601
$ws4.Range("F13").Value = 299
$ws4.Range("F14").Value = 1067
